$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# The "Contact" / "No display for ContactDetail" row becomes
# "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The duplicate "Contact" / "No display for ContactDetail" row is removed
# entirely (rows below shift up one, dimension becomes A1:B14).
$ws.Rows.Item(11).Delete()
